$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Data rows first (this fixes the shared-string table allocation order to
# match the source workbook: data labels before header labels, and the
# "fuel tank" label used by rows 3-5 allocated last).
# ---------------------------------------------------------------------------

# Row 2
$ws.Range("A2").Value = "Geely Borui GE Battery"
$ws.Range("B2").Value = 2000
$ws.Range("C2").Value = "Toyota RAV4 Battery"
$ws.Range("D2").Value = 3500
$ws.Range("E2").Value = "Kia Sportage ICE"
$ws.Range("F2").Value = 2500
$ws.Range("G2").Value = "BMW i8 Motor"
$ws.Range("H2").Value = 1100
$ws.Range("I2").Value = "Master Wheel"
$ws.Range("J2").Value = 100
$ws.Range("K2").Formula = "=SUM(B2,D2,F2,H2,J2)"

# Header row
$ws.Range("A1").Value = "batteryName"
$ws.Range("B1").Value = "batteryCost"
$ws.Range("C1").Value = "fuelTankName"
$ws.Range("D1").Value = "fuelTankCost"
$ws.Range("E1").Value = "ICEName"
$ws.Range("F1").Value = "ICECost"
$ws.Range("G1").Value = "motorName"
$ws.Range("H1").Value = "motorCost"
$ws.Range("I1").Value = "wheelID"
$ws.Range("J1").Value = "wheelCost"
$ws.Range("K1").Value = "total_price"

# Rows 3-5 (same cost structure, different battery/fuel-tank label)
foreach ($r in 3..5) {
    $ws.Range("A$r").Value = "Geely Borui GE Fuel tank"
    $ws.Range("B$r").Value = 500
    $ws.Range("C$r").Value = "Toyota RAV4 Battery"
    $ws.Range("D$r").Value = 3500
    $ws.Range("E$r").Value = "Kia Sportage ICE"
    $ws.Range("F$r").Value = 2500
    $ws.Range("G$r").Value = "BMW i8 Motor"
    $ws.Range("H$r").Value = 1100
    $ws.Range("I$r").Value = "Master Wheel"
    $ws.Range("J$r").Value = 100
    $ws.Range("K$r").Formula = "=SUM(B$r,D$r,F$r,H$r,J$r)"
}

# ---------------------------------------------------------------------------
# Header row styling: Consolas 10, purple font colour, left/middle aligned,
# indented.
# ---------------------------------------------------------------------------
$headerRange = $ws.Range("A1:K1")
$headerRange.Font.Name = "Consolas"
$headerRange.Font.Size = 10
$headerRange.Font.Color = 16058791
$headerRange.HorizontalAlignment = -4131
$headerRange.VerticalAlignment = -4108
$headerRange.IndentLevel = 3

# ---------------------------------------------------------------------------
# Column widths (character units; engine rounds to the nearest 1/6 char, so
# the inputs below are chosen to land as close as possible to the target
# stored widths).
# ---------------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 15.022135416666666
$ws.Columns.Item(2).ColumnWidth = 15.022135416666666
$ws.Columns.Item(3).ColumnWidth = 17.307291666666668
$ws.Columns.Item(4).ColumnWidth = 16.451822916666668
$ws.Columns.Item(5).ColumnWidth = 11.022135416666666
$ws.Columns.Item(6).ColumnWidth = 10.736979166666666
$ws.Columns.Item(7).ColumnWidth = 12.736979166666666
$ws.Columns.Item(8).ColumnWidth = 13.307291666666666
$ws.Columns.Item(9).ColumnWidth = 11.022135416666666
$ws.Columns.Item(10).ColumnWidth = 13.022135416666666
$ws.Columns.Item(11).ColumnWidth = 14.166666666666666

# Selection on K1, matching the saved view state.
[void]$ws.Range("K1").Select()
